# Adding an option for 'all' to country column. Adding country level fuel data.
#
# The "fueldata" sheet gets a new column C ("Country") inserted before the
# existing fuel-type column. Every data row (2-43) gets the value "all" in
# that new column, shifting the old C/D/E columns (fuel type / price /
# emission factor) one column to the right (now D/E/F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fueldata")

# Insert a new column before column C, shifting existing C:E to D:F
# (-4161 == xlShiftToRight)
$ws.Columns.Item(3).Insert(-4161)

# Header for the new column
$ws.Range("C1").Value2 = "Country"

# Fill the new column's data rows with "all"
$ws.Range("C2:C43").Value2 = "all"

# Restore the selection to the newly filled range
$ws.Range("C2:C43").Select()
